$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap M2:R2 <-> M3:R3, and move U3:V3 -> U2:V2, clearing U3:V3
$row2 = $ws.Range("M2:R2").Value2
$row3 = $ws.Range("M3:R3").Value2
$ws.Range("M2:R2").Value = $row3
$ws.Range("M3:R3").Value = $row2

$uv3 = $ws.Range("U3:V3").Value2
$ws.Range("U2:V2").Value = $uv3
$ws.Range("U3:V3").ClearContents()

# View changes: scroll the window so column B is the left-most visible
# column (was J), then move the selection to I11 (was V6).
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("I11").Select() | Out-Null
